$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record above the existing row 162, which pushes
# rows 162..199 down to 163..200 (row 200 ends up holding what used to
# be row 199's data).
$ws.Rows.Item(162).Insert()

$ws.Cells.Item(162, 1).Value = 5
$ws.Cells.Item(162, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(162, 3).Value = "Maule"
$ws.Cells.Item(162, 4).Value = 44511
$ws.Cells.Item(162, 5).Value = 7
$ws.Cells.Item(162, 6).Value = 100112009
$ws.Cells.Item(162, 7).Value = "Acelga"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 500
$ws.Cells.Item(162, 11).Value = 2000
$ws.Cells.Item(162, 12).Value = 2000
$ws.Cells.Item(162, 13).Value = 2000
$ws.Cells.Item(162, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(162, 15).Value = "Región del Maule"
$ws.Cells.Item(162, 16).Value = 500
$ws.Cells.Item(162, 17).Value = 4
$ws.Cells.Item(162, 18).Value = "Hortaliza"
